$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 14 (mailbox 11) now becomes an input from the vector box throttle signal
$ws.Range("D14").Value = "HI_VectorBox"
$ws.Range("C14").Value = "Input"

# Update the active selection left on the sheet after the edit
$ws.Range("D18").Select()
